$wb = $excel.ActiveWorkbook

# --- Sheet "optimal models": add new row for BFT (Pacific bluefin tuna) ---
$ws1 = $wb.Worksheets.Item("optimal models")
$ws1.Range("A9").Value = "BFT"
$ws1.Range("B9").Value = 2
$ws1.Range("C9").Value = 0.5
$ws1.Range("D9").Value = 0.01
$ws1.Range("E9").Value = 3050
$ws1.Range("F9").Value = 0.9982
$ws1.Range("G9").Value = 0.9924
$ws1.Range("H9").Value = 0.9908
$ws1.Range("I9").Value = 0.028539740000000001

# --- Sheet "relative importance": add new row for BFT (Pacific bluefin tuna) ---
$ws2 = $wb.Worksheets.Item("relative importance")
$ws2.Range("A9").Value = "BFT"
$ws2.Range("B9").Value = 1.2188098000000001
$ws2.Range("C9").Value = 2.8458109
$ws2.Range("D9").Value = 1.7593607
$ws2.Range("E9").Value = 1.8752705000000001
$ws2.Range("F9").Value = 2.6685229000000001
$ws2.Range("G9").Value = 11.7668236
$ws2.Range("H9").Value = 1.2960031999999999
$ws2.Range("I9").Value = 1.5037383
$ws2.Range("J9").Value = 0.4955521
$ws2.Range("K9").Value = 6.1918363999999997
$ws2.Range("L9").Value = 19.320694599999999
$ws2.Range("M9").Value = 14.5793008
$ws2.Range("N9").Value = 9.3146749999999994
$ws2.Range("O9").Value = 16.998973899999999
$ws2.Range("P9").Value = 1.3274022999999999
$ws2.Range("Q9").Value = 4.2020749000000004
$ws2.Range("R9").Value = 2.6351502

# --- Update active sheet/tab and selections to match the saved view state ---
# Sheet1 ("optimal models") is no longer the selected tab; its selection moves to A10
$ws1.Range("A10").Select()

# Sheet2 ("relative importance") becomes the active/selected tab; selection at I29
$ws2.Activate()
$ws2.Range("I29").Select()
